$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.85"
$ws.Range("E2").Value = "'1.21%"
$ws.Range("D3").Value = "'41.23"
$ws.Range("E3").Value = "'3.17%"
$ws.Range("D4").Value = "'5.714"
$ws.Range("E4").Value = "'-0.54%"
$ws.Range("D5").Value = "'0.08104"
$ws.Range("E5").Value = "'-0.22%"
$ws.Range("D6").Value = "'8.676"
$ws.Range("E6").Value = "'0.01%"
$ws.Range("D7").Value = "'4.494"
$ws.Range("E7").Value = "'-1.51%"
$ws.Range("D8").Value = "'1.973"
$ws.Range("E8").Value = "'0.40%"
$ws.Range("D9").Value = "'2.999"
$ws.Range("E9").Value = "'0.04%"
$ws.Range("E10").Value = "'-1.36%"
$ws.Range("D11").Value = "'0.1283"
$ws.Range("E11").Value = "'-0.62%"
$ws.Range("D12").Value = "'0.1963"
$ws.Range("E12").Value = "'-1.33%"
$ws.Range("D13").Value = "'8.797"
$ws.Range("E13").Value = "'16.16%"
$ws.Range("D14").Value = "'0.09217"
$ws.Range("E14").Value = "'0.63%"
$ws.Range("D15").Value = "'0.03739"
$ws.Range("E15").Value = "'7.12%"
$ws.Range("D16").Value = "'0.1052"
$ws.Range("E16").Value = "'9.28%"
$ws.Range("D17").Value = "'0.001291"
$ws.Range("E17").Value = "'-2.46%"
$ws.Range("D18").Value = "'0.006330"
$ws.Range("E18").Value = "'2.52%"
$ws.Range("D19").Value = "'3.370"
$ws.Range("E19").Value = "'-0.01%"
$ws.Range("D20").Value = "'0.3500"
$ws.Range("E20").Value = "'-0.92%"
$ws.Range("E21").Value = "'-3.01%"
$ws.Range("D22").Value = "'0.2604"
$ws.Range("E22").Value = "'7.39%"
$ws.Range("D23").Value = "'0.04420"
$ws.Range("E23").Value = "'-0.42%"
$ws.Range("E24").Value = "'-0.01%"
$ws.Range("D25").Value = "'0.004410"
$ws.Range("E25").Value = "'1.09%"
$ws.Range("D26").Value = "'0.0001237"
$ws.Range("E26").Value = "'3.98%"
$ws.Range("D39").Value = "'0.02792"
$ws.Range("E39").Value = "'10.77%"
$ws.Range("D40").Value = "'0.05578"
$ws.Range("E40").Value = "'6.81%"
$ws.Range("D41").Value = "'0.007524"
$ws.Range("E41").Value = "'2.95%"
$ws.Range("D42").Value = "'0.009833"
$ws.Range("E42").Value = "'10.71%"
$ws.Range("D43").Value = "'0.1421"
$ws.Range("E43").Value = "'-0.61%"
$ws.Range("D44").Value = "'0.002102"
$ws.Range("E44").Value = "'-3.99%"
$ws.Range("D45").Value = "'0.01180"
$ws.Range("E45").Value = "'18.11%"
$ws.Range("D46").Value = "'0.00006774"
$ws.Range("E46").Value = "'0.83%"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("D48").Value = "'0.003065"
$ws.Range("E48").Value = "'6.69%"
$ws.Range("D49").Value = "'0.002274"
$ws.Range("E49").Value = "'26.39%"
$ws.Range("D50").Value = "'0.00002095"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0001995"
$ws.Range("E51").Value = "'-0.19%"
